# Update the "想去人数" (want-to-go count) figures in column F across the
# three sheets that carry event data. This mirrors a re-scrape of the
# source data (gh-pages output regeneration) - only column F values change.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$s1 = @{
    2  = 15093
    3  = 19278
    5  = 151
    13 = 61
    14 = 188
    15 = 235
    16 = 69
    17 = 1490
    20 = 102
    21 = 241
    22 = 8074
    24 = 37
    25 = 5
    27 = 1256
    30 = 6094
    31 = 122
    32 = 75
    33 = 175
    34 = 157
    36 = 5501
    37 = 1007
    40 = 53
}
foreach ($row in $s1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $s1[$row]
}

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 21

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$s4 = @{
    2  = 15093
    3  = 19278
    5  = 151
    13 = 61
    14 = 188
    15 = 235
    16 = 69
    17 = 1490
    21 = 102
    22 = 241
    23 = 8074
    25 = 37
    26 = 5
    28 = 1256
    31 = 21
    33 = 6094
    34 = 122
    35 = 75
    36 = 175
    37 = 157
    39 = 5501
    40 = 1007
    43 = 53
}
foreach ($row in $s4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $s4[$row]
}
